$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline-variables")
$ws.Activate()

# Rename existing DNS variable names to their Treasury-prefixed versions
$ws.Range("A23").Value = "tdns1"
$ws.Range("A24").Value = "tdns2"
$ws.Range("A25").Value = "tdns3"

# Row 23: Treasury-FFR Spread Level (10-Year Level)
$ws.Range("B23").Value = "Treasury-FFR Spread Level (10-Year Level)"
$ws.Range("C23").Value = "Interest Rates"
$ws.Range("D23").Value = "calc"
$ws.Range("F23").Value = "%"
$ws.Range("G23").Value = "m"
$ws.Range("I23").Value = "d"
$ws.Range("J23").Value = "base"
$ws.Range("K23").Value = "none"

# Row 24: Treasury-FFR Spread Slope (Negative of 10Y-3M Spread)
$ws.Range("B24").Value = "Treasury-FFR Spread Slope (Negative of 10Y-3M Spread)"
$ws.Range("C24").Value = "Interest Rates"
$ws.Range("D24").Value = "calc"
$ws.Range("F24").Value = "%"
$ws.Range("G24").Value = "m"
$ws.Range("I24").Value = "d"
$ws.Range("J24").Value = "base"
$ws.Range("K24").Value = "none"

# Row 25: Treasury-FFR Spread Curvature
$ws.Range("B25").Value = "Treasury-FFR Spread Curvature"
$ws.Range("C25").Value = "Interest Rates"
$ws.Range("D25").Value = "calc"
$ws.Range("F25").Value = "%"
$ws.Range("G25").Value = "m"
$ws.Range("I25").Value = "d"
$ws.Range("J25").Value = "base"
$ws.Range("K25").Value = "none"

# Update selection to match the authored state
$ws.Range("I4").Select()
